$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row labels: *_old -> *_FV2210, *_new -> *_FV2304 ---
$headerRange = $ws.Range("A1:U1")
for ($i = 1; $i -le $headerRange.Cells.Count; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $val = [string]$cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2210"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
        }
    }
}

# --- Freeze the top row (row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# --- Convert the data range into an Excel Table (ListObject) ---
$tableRange = $ws.Range("A1:U57")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
